$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing header cell (AC1) onto the three
# new header cells so they pick up the same bold/border/center style (s="1").
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value = 94
    $ws.Cells.Item($row, 31).Value = 68
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "done"
